$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "About" sheet: recalibrate the capacity supply curve parameters.
#   - max share            (B12): 0.35 -> 0.3
#   - max profitability    (B16): 15   -> 20
# All of the dependent formulas in B21:B121 recompute automatically.
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B12").Value = 0.3
$wsAbout.Range("B16").Value = 20

# ---------------------------------------------------------------------------
# "CSC-CSCCCMvSoECBtY" sheet: extend the capacity-cost-multiplier curve from
# 15 years out to 20 years out. Column BK used to hold the hard-coded "far
# future" sentinel (1000) directly after BJ (10-year mark, value 15); now the
# 0.25-year step pattern continues through column CD (value 20), and the
# sentinel moves out to column CE.
# ---------------------------------------------------------------------------
$wsCurve = $wb.Worksheets.Item("CSC-CSCCCMvSoECBtY")

$cols = @("BK","BL","BM","BN","BO","BP","BQ","BR","BS","BT","BU","BV","BW","BX","BY","BZ","CA","CB","CC","CD")
$prev = "BJ"
foreach ($col in $cols) {
    # Row 2: running "time" axis, each cell is prior cell + 0.25 years.
    $wsCurve.Range($col + "2").Formula = "=" + $prev + "2+0.25"
    # Row 1: capacity cost multiplier at that time, same formula pattern
    # used by every other column on the row.
    $wsCurve.Range($col + "1").Formula = "=(1-EXP(-((" + $col + "2/About!`$B`$16-(About!`$B`$15-0.5))/About!`$B`$13)^About!`$B`$14))*About!`$B`$12"
    $prev = $col
}

# New terminal "far future" sentinel, now at CE instead of BK.
$wsCurve.Range("CE2").Value = 1000
$wsCurve.Range("CE1").Formula = "=(1-EXP(-((CE2/About!`$B`$16-(About!`$B`$15-0.5))/About!`$B`$13)^About!`$B`$14))*About!`$B`$12"

# ---------------------------------------------------------------------------
# Window/selection state captured by the author when they saved the file.
# Select on CSC-CSCCCMvSoECBtY first, then return to About last so About
# remains the active/tabSelected sheet, matching the saved workbook state.
# ---------------------------------------------------------------------------
$wsCurve.Activate()
$wsCurve.Range("H43").Select()

$wsAbout.Activate()
$wsAbout.Range("H39").Select()
